# Apply crypto price/volume updates per diff (Thu Jan 18 09:11:24 UTC 2024 refresh)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "42.743.67"
$ws.Range("E2").Value = "  -0.12%  "

$ws.Range("D3").Value = "2.538.02"
$ws.Range("E3").Value = "  -0.49%  "

$ws.Range("E4").Value = "  -0.01%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "311.35"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +0.33%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "100.34"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +2.04%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.567"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  -1.00%  "

$ws.Range("E8").Value = "  +0.05%  "

$ws.Range("E9").Value = "  -1.52%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "35.69"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +0.00%  "

$ws.Range("E11").Value = "  -0.28%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "7.31"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -1.32%  "

$ws.Range("E13").Value = "  +1.06%  "

$ws.Range("D14").Value = "2.927.46"
$ws.Range("E14").Value = "  -0.31%  "

$ws.Range("B15").Value = "Chainlink"
$ws.Range("C15").Value = "https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "15.36"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -3.27%  "

$ws.Range("B16").Value = "WrappedEther"
$ws.Range("C16").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D16").Value = "2.551.60"
$ws.Range("E16").Value = "  -2.14%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.816"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -2.83%  "

$ws.Range("D18").Value = "42.727.71"
$ws.Range("E18").Value = "  -0.20%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "6.71"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -0.75%  "

$ws.Range("E20").Value = "  -0.55%  "

$ws.Range("E21").Value = "  -0.42%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "70.10"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +1.05%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "243.39"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -1.86%  "

$ws.Range("E24").Value = "  -1.38%  "

$ws.Range("E25").Value = "  -1.46%  "

$ws.Range("E26").Value = "  +0.05%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "25.61"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -4.17%  "

$ws.Range("E28").Value = "  -1.42%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "10.17"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +0.37%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "38.56"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -4.47%  "

$ws.Range("B31").Value = "Filecoin"
$ws.Range("C31").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "5.87"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +1.94%  "

$ws.Range("B32").Value = "Monero"
$ws.Range("C32").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "158.04"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +0.10%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "2.75"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +5.76%  "

$ws.Range("E34").Value = "  +2.04%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.0793"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -0.99%  "

$ws.Range("B36").Value = "LidoDAOToken"
$ws.Range("C36").Value = "https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "3.15"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -4.48%  "

$ws.Range("B37").Value = "Celestia"
$ws.Range("C37").Value = "https://coinranking.com/coin/YQcD0lBl7+celestia-tia"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "18.09"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -1.82%  "

$ws.Range("E38").Value = "  -5.68%  "

$ws.Range("E39").Value = "  -0.26%  "

$ws.Range("E40").Value = "  -0.28%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "4.15"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +1.21%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "21.89"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -2.46%  "

$ws.Range("E43").Value = "  +0.12%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "3.28"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +1.98%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.0298"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -0.47%  "

$ws.Range("D46").Value = "1.995.06"
$ws.Range("E46").Value = "  -0.31%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "9.09"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +0.55%  "

$ws.Range("D48").Value = "2.781.48"

$ws.Range("E49").Value = "  -0.44%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "79.96"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -1.76%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "72.46"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -1.65%  "
